$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: set resultado/profit
$ws.Range("G21").Value = "Fallo"
$ws.Range("H21").Value = -1

# Row 22: set resultado/profit
$ws.Range("G22").Value = "Acierto"
$ws.Range("H22").Value = 2.4

# Row 24: A24 was stored as text, convert to a real number
$ws.Range("A24").Value = 14687083

# Row 25: A25 was stored as text, convert to a real number
$ws.Range("A25").Value = 14721398
